# Updating filtered feeds from workflow
#
# Appends two new feed rows (32 and 33) to the "Filtered Feeds" sheet for the
# new Thermo Fisher / Boehringer Ingelheim NGS-CDx story, one row per source
# (genomeweb.com, 360dx.com), mirroring the existing link/keyword/title
# layout used by every other row in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$genomewebLink = "https://www.genomeweb.com/sequencing/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug"
$dx360Link     = "https://www.360dx.com/sequencing/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug"
$keyword       = "CDx"
$title         = "Thermo Fisher Nabs FDA Approval for NGS-Based CDx for Boehringer Ingelheim Lung Cancer Drug"

# Fill column-by-column (link, then keyword, then title) so the two brand
# new link strings are interned before the new title string -- keeping the
# new title the very last new entry added to the shared-string table, same
# relative order as the rest of the link/keyword/title table.

# --- Column A (link) --------------------------------------------------------
$ws.Range("A32").Value = $genomewebLink
$ws.Range("A33").Value = $dx360Link

# --- Column B (keywords) ----------------------------------------------------
$ws.Range("B32").Value = $keyword
$ws.Range("B33").Value = $keyword

# --- Column C (title) -------------------------------------------------------
$ws.Range("C32").Value = $title
$ws.Range("C33").Value = $title

# Wire up the link-column hyperlinks (values already match the target URL,
# same convention as the existing rows) and register the relationships.
$ws.Hyperlinks.Add($ws.Range("A32"), $genomewebLink)
$ws.Hyperlinks.Add($ws.Range("A33"), $dx360Link)

# Match the hyperlink-cell style used by the rest of column A.
$ws.Range("A32").Style = $ws.Range("A31").Style
$ws.Range("A33").Style = $ws.Range("A31").Style
